$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update Approved/Rejected column: "Reworked" -> "Approved" for the two test scenarios
$ws.Range("I6").Value = "Approved"
$ws.Range("I12").Value = "Approved"

# Update the active selection to match the recorded cursor position
$ws.Range("I12").Select()
